$d = $word.ActiveDocument
# 1. Remove centering from the empty paragraph following the title
$rng = $d.Content
$found = $rng.Find.Execute("Wood Budget Analysis - Progress Report", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titlePara = $rng.Paragraphs(1)
$para2 = $titlePara.Next()
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="332935AB" w14:textId="4EFCB5DD" w:rsidR="004745FC" w:rsidRDefault="004745FC" w:rsidP="004745FC"><w:pPr><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$para2.Range.InsertXML($xml)

# 2. Remove lastRenderedPageBreak before "* General questions about mortality"
$rng = $d.Content
$found = $rng.Find.Execute("General questions about mortality", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gqPara = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="7E0F2D56" w14:textId="4ED6C918" w:rsidR="004C7CBF" w:rsidRDefault="004C7CBF" w:rsidP="004C7CBF"><w:pPr><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="004C7CBF"><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">* General questions about mortality </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>que</w:t></w:r><w:r w:rsidRPr="004C7CBF"><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>tions</w:t></w:r><w:r w:rsidRPr="004C7CBF"><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>…</w:t></w:r></w:p>'
$gqPara.Range.InsertXML($xml)

# 3. Add lastRenderedPageBreak before "slopes.kmz"
$rng = $d.Content
$found = $rng.Find.Execute("slopes.kmz", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$slopesPara = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="2620510D" w14:textId="3087ECFE" w:rsidR="00D01DB6" w:rsidRDefault="00D01DB6" w:rsidP="004C7CBF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00D01DB6"><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>slopes.kmz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>kmz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file showing the top and bottom of the transects used in google earth</w:t></w:r></w:p>'
$slopesPara.Range.InsertXML($xml)

# 4. Remove lastRenderedPageBreak before "Cutblock_species.csv"
$rng = $d.Content
$found = $rng.Find.Execute("Cutblock_species.csv", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cutblockPara = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="699B9819" w14:textId="197D821D" w:rsidR="00D01DB6" w:rsidRDefault="00D01DB6" w:rsidP="004C7CBF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00D01DB6"><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Cutblock_species.csv</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – excel file showing the characteristics of the polygon used to derive estimates of the forest conditions for the harvested </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>cutblock</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> for the tipsy model</w:t></w:r></w:p>'
$cutblockPara.Range.InsertXML($xml)

# 5. Remove the _GoBack bookmark from its old location (it is re-added elsewhere)
$rng = $d.Content
$found = $rng.Find.Execute("set to 0.004, the estimate for lodgepole pine", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPara = $rng.Paragraphs(1).Next()
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="50F81551" w14:textId="77777777" w:rsidR="008D1286" w:rsidRPr="00E86EEE" w:rsidRDefault="008D1286" w:rsidP="00E86EEE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>'
$bmPara.Range.InsertXML($xml)

# 6. Add lastRenderedPageBreak before "I" in the Hassan-paper paragraph
$rng = $d.Content
$found = $rng.Find.Execute("In the Hassan paper, D", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hassanPara = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="7B82E0B4" w14:textId="3B56B734" w:rsidR="00324BF1" w:rsidRPr="004C7CBF" w:rsidRDefault="004C7CBF" w:rsidP="00324BF1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>I</w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>n the Hassan paper, D</w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">it is referred to as density of trees </w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>(with units of m</w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="00324BF1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>/ha). I think it should be a volume</w:t></w:r><w:r w:rsidR="00E86EEE"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>…therefore the m3/ha data was used in code and converted to m3/m2</w:t></w:r></w:p>'
$hassanPara.Range.InsertXML($xml)

# 7. Replace the "Update on wood budget code:" paragraph with the new content
$rng = $d.Content
$found = $rng.Find.Execute("Update on wood budget code:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Two scripts are required to calculate the wood budget. The script “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>input_functions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">” contains all the required functions to run the model. The script “wood-budget” loads the required data and runs the functions to calculate the input for </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">budget. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Update on wood budget </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>analysis</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p>'
$targetPara.Range.InsertXML($xml)

Write-Output "done"
